$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data rows (keep header row 1 intact)
$ws.Range("A2:C18").ClearContents()

$data = @(
  @("Dyson Daniels","PG,SG,SF","Atlanta Hawks"),
  @("Kelly Oubre Jr.","SG,SF","Philadelphia 76ers"),
  @("Devin Vassell","SG,SF","San Antonio Spurs"),
  @("Jaden McDaniels","SF,PF","Minnesota Timberwolves"),
  @("Toumani Camara","SG,SF,PF","Portland Trail Blazers"),
  @("Michael Porter Jr.","SF,PF","Denver Nuggets"),
  @("Bam Adebayo","PF,C","Miami Heat"),
  @("Guerschon Yabusele","PF,C","Philadelphia 76ers"),
  @("Kel'el Ware","PF,C","Miami Heat"),
  @("Josh Hart","SG,SF,PF","New York Knicks"),
  @("Donovan Mitchell","PG,SG","Cleveland Cavaliers"),
  @("Malik Beasley","SG,SF","Detroit Pistons"),
  @("De'Andre Hunter","SF,PF","Cleveland Cavaliers"),
  @("Julius Randle","PF,C","Minnesota Timberwolves"),
  @("Cam Thomas","SG,SF","Brooklyn Nets"),
  @("Domantas Sabonis","C","Sacramento Kings"),
  @("Kristaps Porzingis","PF,C","Boston Celtics"),
  @("Lauri Markkanen","SF,PF","Utah Jazz")
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

Write-Output "done"
